$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update beta_distance_range Max value (C3): 8.9 -> 10.1
$ws.Range("C3").Value = 10.1

# Delete the theta_threshold_range row (row 5), shifting pie_threshold_range up to row 5
$ws.Rows("5").Delete()

# Update new row 5 (pie_threshold_range) Min value (B5): 2 -> 0
$ws.Range("B5").Value = 0

# Update selection to C4
$ws.Range("C4").Select()

# Set up page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
